$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.234.76"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.317.01"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'541.63"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'132.72"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "2.314.70"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'23.99"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "2.728.73"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "59.086.68"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "2.315.43"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").Value = "'313.47"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "'6.60"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'62.67"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "'1.73"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  +6.11%  "
$ws.Range("D31").Value = "'170.36"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").Value = "0.0₃0743"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "'5.90"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "'0.386"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  +6.04%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'17.87"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  +3.16%  "
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").Value = "'305.76"
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("D42").Value = "'141.19"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").Value = "'0.0960"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").Value = "'0.0497"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'0.0213"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").Value = "'11.00"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").Value = "'4.64"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  +2.09%  "
